$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "oneday"
$ws.Range("B2").Value = "https://www.google.com/aclk?sa=l&ai=DChcSEwi8nNbXgv2CAxWykmgJHTOCDG0YABAAGgJ3Zg&ase=2&gclid=EAIaIQobChMIvJzW14L9ggMVspJoCR0zggxtEBAYASAAEgJAHvD_BwE&sig=AOD64_2E9U-QK_qGM4hW0cWuw4HufTXQRQ&adurl&ctype=99"
$ws.Range("A3").Value = "makeagency.co"
$ws.Range("A4").Value = "sidedishmedia.co"
$ws.Range("A6").Value = "bathmarketingcompany.uixweb"
$ws.Range("A7").Value = "bbh-usa"
$ws.Range("A8").Value = "generationmedia.co"
$ws.Range("A9").Value = "saatchi.co"
$ws.Range("A12").Value = "ality.co"
$ws.Range("A16").Value = "tipigroup"
